$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.778.31'
$ws.Range('E2').Value = '  -3.33%  '
$ws.Range('D3').Value = '2.912.68'
$ws.Range('E3').Value = '  -3.94%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.07'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.88'
$ws.Range('E6').Value = '  -6.44%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -2.77%  '
$ws.Range('D9').Value = '2.912.38'
$ws.Range('E9').Value = '  -3.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.80'
$ws.Range('E10').Value = '  +5.91%  '
$ws.Range('E11').Value = '  -4.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.446'
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.35'
$ws.Range('E14').Value = '  -6.29%  '
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').Value = '3.396.19'
$ws.Range('E16').Value = '  -3.98%  '
$ws.Range('D17').Value = '60.752.42'
$ws.Range('E18').Value = '  -5.16%  '
$ws.Range('D19').Value = '2.911.76'
$ws.Range('E19').Value = '  -3.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '430.28'
$ws.Range('E20').Value = '  -4.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.62'
$ws.Range('E21').Value = '  -4.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.683'
$ws.Range('E22').Value = '  -1.84%  '
$ws.Range('E23').Value = '  -4.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.33'
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.84'
$ws.Range('E25').Value = '  -3.07%  '
$ws.Range('E26').Value = '  -4.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.86'
$ws.Range('E27').Value = '  -3.98%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.61'
$ws.Range('E30').Value = '  -3.00%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.16'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('E33').Value = '  -3.83%  '
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('D35').Value = '0.0₃0868'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.00'
$ws.Range('E38').Value = '  -5.82%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.126'
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.71'
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('E41').Value = '  -5.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.66'
$ws.Range('E42').Value = '  -4.61%  '
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('E45').Value = '  -3.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '374.85'
$ws.Range('E46').Value = '  -5.35%  '
$ws.Range('D47').Value = '2.668.85'
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.73'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.27'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('E51').Value = '  -1.94%  '
